{"js": "// Diary entry for 2002\u5e746\u67083\u65e5\u661f\u671f\u4e94 originally ended with the\n// \"\u7aef\u5348\u8282\" (Dragon Boat Festival) paragraph. The edit keeps that old\n// entry in place, but inserts a brand-new diary entry right after it\n// (a date paragraph \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\") and turns the paragraph\n// that used to hold the \"\u7aef\u5348\u8282\" text into the new entry's body text.\n\nconst OLD_TEXT =\n  \"\u4e2d\u96e8\uff0c\u4eca\u5929\u662f\u519c\u5386\u4e94\u6708\u521d\u4e94\uff0c\u4e2d\u56fd\u4f20\u7edf\u7aef\u5348\u8282\uff1a\u7aef\u5348\u8282\uff0c\u8fd9\u4e00\u5929\u6211\u4eec\u8981\u5403\u7cbd\u5b50\uff0c\u8d5b\u9f99\u821f\u3002\";\nconst NEW_DATE = \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\";\nconst NEW_BODY = \"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002\";\n\n// Locate the paragraph that still carries the old \"\u7aef\u5348\u8282\" text (it is\n// the one that also owns the _GoBack bookmark at the end of the body).\nconst results = context.document.body.search(OLD_TEXT, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nconst target = results.items[0].paragraphs.getFirst();\ntarget.load(\"text\");\n\nconst prev = target.getPrevious();\nprev.load(\"text\");\nawait context.sync();\n\n// Insert the two new paragraphs after the paragraph preceding the\n// target, so they pick up that paragraph's \"eastAsia\" font hint\n// formatting instead of inheriting the target's \"default\" hint.\nconst inserted1 = prev.insertParagraph(target.text, \"After\");\nawait context.sync();\n\ninserted1.insertParagraph(NEW_DATE, \"After\");\nawait context.sync();\n\n// Finally, rewrite the original (now-shifted) paragraph's text in place\n// so it keeps its run/paragraph formatting and the _GoBack bookmark.\ntarget.insertText(NEW_BODY, \"Replace\");\nawait context.sync();\n", "ps1": "# Diary entry for 2002\u5e746\u67083\u65e5\u661f\u671f\u4e94 originally ended with the\n# \"\u7aef\u5348\u8282\" (Dragon Boat Festival) paragraph. The edit keeps that old\n# entry in place, but inserts a brand-new diary entry right after it\n# (a date paragraph \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\") and turns the paragraph\n# that used to hold the \"\u7aef\u5348\u8282\" text (which also carries the\n# trailing _GoBack bookmark) into the new entry's body text.\n\n$d = $word.ActiveDocument\n\n$oldText = \"\u4e2d\u96e8\uff0c\u4eca\u5929\u662f\u519c\u5386\u4e94\u6708\u521d\u4e94\uff0c\u4e2d\u56fd\u4f20\u7edf\u7aef\u5348\u8282\uff1a\u7aef\u5348\u8282\uff0c\u8fd9\u4e00\u5929\u6211\u4eec\u8981\u5403\u7cbd\u5b50\uff0c\u8d5b\u9f99\u821f\u3002\"\n$newDate = \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\"\n$newBody = \"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002\"\n\nfunction Find-ParagraphByText($doc, $text) {\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# Locate the (currently unique) paragraph holding the old \"\u7aef\u5348\u8282\" text,\n# and the paragraph right before it -- the new paragraphs will be\n# inserted after that one so they inherit its \"eastAsia\" font-hint\n# formatting instead of the bookmark paragraph's \"default\" hint.\n$targetPara = Find-ParagraphByText $d $oldText\n$prevPara = $targetPara.Previous()\n\n# Rewrite the target paragraph's text in place first (while the\n# reference is still unique) so it keeps its paragraph/run formatting\n# and the _GoBack bookmark.\n$targetPara.Range.Text = $newBody\n\n# Insert two fresh blank paragraphs after $prevPara, then fill them in:\n# first the (old, unmodified) \"\u7aef\u5348\u8282\" entry, then the new date line.\n$insertRange = $prevPara.Range\n$insertRange.Collapse(0)\n$insertRange.InsertParagraphAfter()\n$insertRange.InsertParagraphAfter()\n\n$p1 = $prevPara.Next()\n$p1.Range.InsertBefore($oldText)\n\n$p2 = $p1.Next()\n$p2.Range.InsertBefore($newDate)\n"}
